$wb = $excel.ActiveWorkbook

$wsImport = $wb.Worksheets.Item("Import Priorities")
$wsCoupling = $wb.Worksheets.Item("Coupling Parameters")

# Add the new "start_year_fuel_trends" config entry as row 9 of the
# "Coupling Parameters" sheet.
$wsCoupling.Range("A9").Value = "start_year_fuel_trends"
$wsCoupling.Range("B9").Value = 5

# Update selections to reflect where the editor left off, and make the
# "Coupling Parameters" sheet the active tab.
[void]$wsImport.Range("A11:B11").Select()

[void]$wsCoupling.Activate()
[void]$wsCoupling.Range("B9").Select()
